# Agregue reporte de egresos e ingresos
# Fill in the "Estado" column (D) for the remaining requirement rows and
# color-code columns A:D to match each row's status (green = DESARROLLADO,
# yellow = EN PROCESO, no fill = SIN INICIAR).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$verde    = 5287936   # RGB(0,176,80)  -> DESARROLLADO
$amarillo = 65535     # RGB(255,255,0) -> EN PROCESO

# Rows finished (DESARROLLADO) and highlighted in green across A:D
$filasDesarrollado = @(4, 5, 6, 7, 8, 9, 17)
foreach ($r in $filasDesarrollado) {
    $rango = $ws.Range("A" + $r + ":D" + $r)
    $rango.Interior.Color = $verde
    $ws.Range("D" + $r).Value = "DESARROLLADO"
}

# Row 10 is already finished but keeps its original (no-fill) formatting
$ws.Range("D10").Value = "DESARROLLADO"

# Rows in progress (EN PROCESO) and highlighted in yellow across A:D
$filasEnProceso = @(11, 12)
foreach ($r in $filasEnProceso) {
    $rango = $ws.Range("A" + $r + ":D" + $r)
    $rango.Interior.Color = $amarillo
    $ws.Range("D" + $r).Value = "EN PROCESO"
}

# Rows not started yet (SIN INICIAR) keep their original (no-fill) formatting
$filasSinIniciar = @(13, 14, 15, 16)
foreach ($r in $filasSinIniciar) {
    $ws.Range("D" + $r).Value = "SIN INICIAR"
}

# Restore the last active selection recorded for the sheet
$ws.Range("B22").Select()
